$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110, shifting existing rows 110..142 down to 111..143.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new data record.
$ws.Cells.Item(110, 1).Value = 5
$ws.Cells.Item(110, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(110, 3).Value = "Maule"
$ws.Cells.Item(110, 4).Value = 44463
$ws.Cells.Item(110, 5).Value = 7
$ws.Cells.Item(110, 6).Value = 100112008
$ws.Cells.Item(110, 7).Value = "Coliflor"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 3000
$ws.Cells.Item(110, 11).Value = 600
$ws.Cells.Item(110, 12).Value = 600
$ws.Cells.Item(110, 13).Value = 600
$ws.Cells.Item(110, 14).Value = "$/unidad"
$ws.Cells.Item(110, 15).Value = "Región del Maule"
$ws.Cells.Item(110, 16).Value = 600
$ws.Cells.Item(110, 17).Value = 1
$ws.Cells.Item(110, 18).Value = "Hortaliza"

Write-Output "Inserted new row 110; sheet now has $($ws.Rows.Count) potential rows (used range updated)."
